$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-like updates (coin names / links) - plain string assignment
$textUpdates = @(
    @{Cell='B8'; Value='GateToken'},
    @{Cell='C8'; Value='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'},
    @{Cell='B9'; Value='MXToken'},
    @{Cell='C9'; Value='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Cell='B10'; Value='LiechtensteinCryptoassetsExchange'},
    @{Cell='C10'; Value='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'},
    @{Cell='B11'; Value='WazirX'},
    @{Cell='C11'; Value='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'},
    @{Cell='B12'; Value='MandalaExchangeToken'},
    @{Cell='C12'; Value='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'},
    @{Cell='B13'; Value='BitrueCoin'},
    @{Cell='C13'; Value='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'},
    @{Cell='B14'; Value='BitMartToken'},
    @{Cell='C14'; Value='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'},
    @{Cell='B15'; Value='BitForexToken'},
    @{Cell='C15'; Value='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'},
    @{Cell='B16'; Value='TigerCash'},
    @{Cell='C16'; Value='https://coinranking.com/coin/6hIn06L2+tigercash-tch'},
    @{Cell='B17'; Value='LEO'},
    @{Cell='C17'; Value='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'}
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Numeric-looking updates (price / volume percentages) - force text storage
# so the cell keeps its original "plain text" representation (e.g. "302.92",
# "2.85%") instead of being auto-converted to a Number/Percentage by Excel.
$numericUpdates = @(
    @{Cell='D2'; Value='302.92'},
    @{Cell='E2'; Value='2.85%'},
    @{Cell='D3'; Value='33.96'},
    @{Cell='E3'; Value='9.36%'},
    @{Cell='D4'; Value='5.082'},
    @{Cell='E4'; Value='2.98%'},
    @{Cell='D5'; Value='0.07815'},
    @{Cell='E5'; Value='6.20%'},
    @{Cell='D6'; Value='2.450'},
    @{Cell='E6'; Value='5.89%'},
    @{Cell='D7'; Value='8.004'},
    @{Cell='E7'; Value='4.01%'},
    @{Cell='D8'; Value='3.898'},
    @{Cell='E8'; Value='3.67%'},
    @{Cell='D9'; Value='0.9384'},
    @{Cell='E9'; Value='2.76%'},
    @{Cell='D10'; Value='0.09950'},
    @{Cell='E10'; Value='21.06%'},
    @{Cell='D11'; Value='0.1778'},
    @{Cell='E11'; Value='5.43%'},
    @{Cell='D12'; Value='0.08559'},
    @{Cell='E12'; Value='4.15%'},
    @{Cell='D13'; Value='0.03342'},
    @{Cell='E13'; Value='6.56%'},
    @{Cell='D14'; Value='0.09919'},
    @{Cell='E14'; Value='-1.54%'},
    @{Cell='D15'; Value='0.001476'},
    @{Cell='E15'; Value='-1.99%'},
    @{Cell='D16'; Value='0.005686'},
    @{Cell='E16'; Value='-0.77%'},
    @{Cell='D17'; Value='3.465'},
    @{Cell='E17'; Value='-0.50%'},
    @{Cell='E18'; Value='5.33%'},
    @{Cell='D19'; Value='0.3356'},
    @{Cell='E19'; Value='0.79%'},
    @{Cell='D20'; Value='0.1345'},
    @{Cell='E20'; Value='3.16%'},
    @{Cell='D21'; Value='4.303'},
    @{Cell='E22'; Value='-0.51%'},
    @{Cell='D23'; Value='0.04627'},
    @{Cell='E23'; Value='1.66%'},
    @{Cell='D24'; Value='0.001220'},
    @{Cell='E24'; Value='0.69%'},
    @{Cell='D25'; Value='0.004410'},
    @{Cell='E25'; Value='1.64%'},
    @{Cell='D26'; Value='0.0001297'},
    @{Cell='E26'; Value='-0.26%'},
    @{Cell='D27'; Value='0.0003377'},
    @{Cell='E27'; Value='-0.55%'},
    @{Cell='E39'; Value='7.53%'},
    @{Cell='D40'; Value='0.04826'},
    @{Cell='E40'; Value='8.48%'},
    @{Cell='D41'; Value='0.007831'},
    @{Cell='E41'; Value='6.40%'},
    @{Cell='D42'; Value='0.009783'},
    @{Cell='E42'; Value='10.83%'},
    @{Cell='D43'; Value='0.1412'},
    @{Cell='E43'; Value='6.60%'},
    @{Cell='D44'; Value='0.002074'},
    @{Cell='E44'; Value='-0.35%'},
    @{Cell='D45'; Value='0.009121'},
    @{Cell='E45'; Value='0.14%'},
    @{Cell='D46'; Value='0.00006079'},
    @{Cell='E46'; Value='3.18%'},
    @{Cell='D47'; Value='0.00000000748'},
    @{Cell='E47'; Value='-0.27%'},
    @{Cell='D48'; Value='2.794'},
    @{Cell='E48'; Value='24.67%'},
    @{Cell='E49'; Value='-31.36%'},
    @{Cell='D50'; Value='0.00002095'},
    @{Cell='E50'; Value='-0.27%'},
    @{Cell='D51'; Value='0.0001996'},
    @{Cell='E51'; Value='-0.27%'}
)

foreach ($u in $numericUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}